$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Value = "2021-07-23"
$ws.Cells.Item(2,11).Value = "Start Ruby"
$ws.Cells.Item(2,12).Value = "Primera"
$ws.Cells.Item(2,13).Value = 140
$ws.Cells.Item(2,14).Value = 9800
$ws.Cells.Item(2,15).Value = 9800
$ws.Cells.Item(2,16).Value = 9800
$ws.Cells.Item(2,17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(2,18).Value = "Región de O'Higgins"
$ws.Cells.Item(2,19).Value = 700
$ws.Cells.Item(2,20).Value = 14

# Row 3
$ws.Cells.Item(3,4).Value = "2020-12-24"
$ws.Cells.Item(3,11).Value = "Start Ruby"
$ws.Cells.Item(3,12).Value = "Primera"
$ws.Cells.Item(3,13).Value = 16
$ws.Cells.Item(3,14).Value = 150000
$ws.Cells.Item(3,15).Value = 150000
$ws.Cells.Item(3,16).Value = 150000
$ws.Cells.Item(3,17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(3,18).Value = "Provincia de Limarí"
$ws.Cells.Item(3,19).Value = 429
$ws.Cells.Item(3,20).Value = 350

# Row 4
$ws.Cells.Item(4,4).Value = "2021-01-12"
$ws.Cells.Item(4,11).Value = "Start Ruby"
$ws.Cells.Item(4,12).Value = "Primera"
$ws.Cells.Item(4,13).Value = 16
$ws.Cells.Item(4,14).Value = 180000
$ws.Cells.Item(4,15).Value = 180000
$ws.Cells.Item(4,16).Value = 180000
$ws.Cells.Item(4,17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(4,18).Value = "Región Metropolitana"
$ws.Cells.Item(4,19).Value = 514
$ws.Cells.Item(4,20).Value = 350

# Row 5
$ws.Cells.Item(5,4).Value = "2021-04-23"
$ws.Cells.Item(5,11).Value = "Start Ruby"
$ws.Cells.Item(5,12).Value = "Primera"
$ws.Cells.Item(5,13).Value = 16
$ws.Cells.Item(5,14).Value = 350000
$ws.Cells.Item(5,15).Value = 350000
$ws.Cells.Item(5,16).Value = 350000
$ws.Cells.Item(5,17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(5,18).Value = "Región Metropolitana"
$ws.Cells.Item(5,19).Value = 1000
$ws.Cells.Item(5,20).Value = 350

# Row 6
$ws.Cells.Item(6,4).Value = "2020-12-31"
$ws.Cells.Item(6,11).Value = "Red Blush"
$ws.Cells.Item(6,12).Value = "Primera"
$ws.Cells.Item(6,13).Value = 12
$ws.Cells.Item(6,14).Value = 130000
$ws.Cells.Item(6,15).Value = 130000
$ws.Cells.Item(6,16).Value = 130000
$ws.Cells.Item(6,17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(6,18).Value = "Provincia de Limarí"
$ws.Cells.Item(6,19).Value = 371
$ws.Cells.Item(6,20).Value = 350

# Row 7
$ws.Cells.Item(7,4).Value = "2020-12-02"
$ws.Cells.Item(7,11).Value = "Start Ruby"
$ws.Cells.Item(7,12).Value = "Primera"
$ws.Cells.Item(7,13).Value = 140
$ws.Cells.Item(7,14).Value = 9800
$ws.Cells.Item(7,15).Value = 9800
$ws.Cells.Item(7,16).Value = 9800
$ws.Cells.Item(7,17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(7,18).Value = "Región de O'Higgins"
$ws.Cells.Item(7,19).Value = 700
$ws.Cells.Item(7,20).Value = 14

# Row 8
$ws.Cells.Item(8,4).Value = "2021-06-29"
$ws.Cells.Item(8,11).Value = "Start Ruby"
$ws.Cells.Item(8,12).Value = "Primera"
$ws.Cells.Item(8,13).Value = 20
$ws.Cells.Item(8,14).Value = 180000
$ws.Cells.Item(8,15).Value = 180000
$ws.Cells.Item(8,16).Value = 180000
$ws.Cells.Item(8,17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(8,18).Value = "Hijuelas"
$ws.Cells.Item(8,19).Value = 514
$ws.Cells.Item(8,20).Value = 350

# Row 9
$ws.Cells.Item(9,4).Value = "2021-06-29"
$ws.Cells.Item(9,11).Value = "Start Ruby"
$ws.Cells.Item(9,12).Value = "Segunda"
$ws.Cells.Item(9,13).Value = 16
$ws.Cells.Item(9,14).Value = 150000
$ws.Cells.Item(9,15).Value = 150000
$ws.Cells.Item(9,16).Value = 150000
$ws.Cells.Item(9,17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(9,18).Value = "Provincia de Limarí"
$ws.Cells.Item(9,19).Value = 429
$ws.Cells.Item(9,20).Value = 350

# Row 10
$ws.Cells.Item(10,4).Value = "2021-09-07"
$ws.Cells.Item(10,11).Value = "Start Ruby"
$ws.Cells.Item(10,12).Value = "Primera"
$ws.Cells.Item(10,13).Value = 14
$ws.Cells.Item(10,14).Value = 150000
$ws.Cells.Item(10,15).Value = 160000
$ws.Cells.Item(10,16).Value = 155000
$ws.Cells.Item(10,17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(10,18).Value = "Región de O'Higgins"
$ws.Cells.Item(10,19).Value = 443
$ws.Cells.Item(10,20).Value = 350

# Row 11
$ws.Cells.Item(11,4).Value = "2021-04-22"
$ws.Cells.Item(11,11).Value = "Start Ruby"
$ws.Cells.Item(11,12).Value = "Primera"
$ws.Cells.Item(11,13).Value = 20
$ws.Cells.Item(11,14).Value = 280000
$ws.Cells.Item(11,15).Value = 280000
$ws.Cells.Item(11,16).Value = 280000
$ws.Cells.Item(11,17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(11,18).Value = "Región Metropolitana"
$ws.Cells.Item(11,19).Value = 800
$ws.Cells.Item(11,20).Value = 350

# Row 12
$ws.Cells.Item(12,4).Value = "2020-12-30"
$ws.Cells.Item(12,11).Value = "Start Ruby"
$ws.Cells.Item(12,12).Value = "Primera"
$ws.Cells.Item(12,13).Value = 20
$ws.Cells.Item(12,14).Value = 200000
$ws.Cells.Item(12,15).Value = 210000
$ws.Cells.Item(12,16).Value = 206000
$ws.Cells.Item(12,17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(12,18).Value = "Región de O'Higgins"
$ws.Cells.Item(12,19).Value = 589
$ws.Cells.Item(12,20).Value = 350

# Row 13
$ws.Cells.Item(13,4).Value = "2021-06-09"
$ws.Cells.Item(13,11).Value = "Start Ruby"
$ws.Cells.Item(13,12).Value = "Primera"
$ws.Cells.Item(13,13).Value = 24
$ws.Cells.Item(13,14).Value = 200000
$ws.Cells.Item(13,15).Value = 230000
$ws.Cells.Item(13,16).Value = 215000
$ws.Cells.Item(13,17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(13,18).Value = "Región Metropolitana"
$ws.Cells.Item(13,19).Value = 614
$ws.Cells.Item(13,20).Value = 350

# Row 14
$ws.Cells.Item(14,4).Value = "2021-06-16"
$ws.Cells.Item(14,11).Value = "Start Ruby"
$ws.Cells.Item(14,12).Value = "Primera"
$ws.Cells.Item(14,13).Value = 20
$ws.Cells.Item(14,14).Value = 200000
$ws.Cells.Item(14,15).Value = 230000
$ws.Cells.Item(14,16).Value = 215000
$ws.Cells.Item(14,17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(14,18).Value = "Provincia de Limarí"
$ws.Cells.Item(14,19).Value = 614
$ws.Cells.Item(14,20).Value = 350

# Row 15
$ws.Cells.Item(15,4).Value = "2021-04-12"
$ws.Cells.Item(15,11).Value = "Start Ruby"
$ws.Cells.Item(15,12).Value = "Especial"
$ws.Cells.Item(15,13).Value = 15
$ws.Cells.Item(15,14).Value = 450000
$ws.Cells.Item(15,15).Value = 450000
$ws.Cells.Item(15,16).Value = 450000
$ws.Cells.Item(15,17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(15,18).Value = "Región Metropolitana"
$ws.Cells.Item(15,19).Value = 1286
$ws.Cells.Item(15,20).Value = 350

# Row 16
$ws.Cells.Item(16,4).Value = "2021-04-12"
$ws.Cells.Item(16,11).Value = "Start Ruby"
$ws.Cells.Item(16,12).Value = "Primera"
$ws.Cells.Item(16,13).Value = 20
$ws.Cells.Item(16,14).Value = 430000
$ws.Cells.Item(16,15).Value = 430000
$ws.Cells.Item(16,16).Value = 430000
$ws.Cells.Item(16,17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(16,18).Value = "Región Metropolitana"
$ws.Cells.Item(16,19).Value = 1229
$ws.Cells.Item(16,20).Value = 350

# Row 17
$ws.Cells.Item(17,4).Value = "2021-01-05"
$ws.Cells.Item(17,11).Value = "Start Ruby"
$ws.Cells.Item(17,12).Value = "Especial"
$ws.Cells.Item(17,13).Value = 8
$ws.Cells.Item(17,14).Value = 200000
$ws.Cells.Item(17,15).Value = 200000
$ws.Cells.Item(17,16).Value = 200000
$ws.Cells.Item(17,17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(17,18).Value = "Región de O'Higgins"
$ws.Cells.Item(17,19).Value = 571
$ws.Cells.Item(17,20).Value = 350

# Row 18
$ws.Cells.Item(18,4).Value = "2021-01-05"
$ws.Cells.Item(18,11).Value = "Start Ruby"
$ws.Cells.Item(18,12).Value = "Primera"
$ws.Cells.Item(18,13).Value = 16
$ws.Cells.Item(18,14).Value = 170000
$ws.Cells.Item(18,15).Value = 170000
$ws.Cells.Item(18,16).Value = 170000
$ws.Cells.Item(18,17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(18,18).Value = "Región de O'Higgins"
$ws.Cells.Item(18,19).Value = 486
$ws.Cells.Item(18,20).Value = 350

# Row 19
$ws.Cells.Item(19,4).Value = "2020-12-28"
$ws.Cells.Item(19,11).Value = "Start Ruby"
$ws.Cells.Item(19,12).Value = "Primera"
$ws.Cells.Item(19,13).Value = 8
$ws.Cells.Item(19,14).Value = 150000
$ws.Cells.Item(19,15).Value = 150000
$ws.Cells.Item(19,16).Value = 150000
$ws.Cells.Item(19,17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(19,18).Value = "Región Metropolitana"
$ws.Cells.Item(19,19).Value = 429
$ws.Cells.Item(19,20).Value = 350

# Row 20
$ws.Cells.Item(20,4).Value = "2021-07-12"
$ws.Cells.Item(20,11).Value = "Start Ruby"
$ws.Cells.Item(20,12).Value = "Especial"
$ws.Cells.Item(20,13).Value = 18
$ws.Cells.Item(20,14).Value = 200000
$ws.Cells.Item(20,15).Value = 200000
$ws.Cells.Item(20,16).Value = 200000
$ws.Cells.Item(20,17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(20,18).Value = "Provincia de Quillota"
$ws.Cells.Item(20,19).Value = 571
$ws.Cells.Item(20,20).Value = 350

# Row 21
$ws.Cells.Item(21,4).Value = "2021-04-26"
$ws.Cells.Item(21,11).Value = "Start Ruby"
$ws.Cells.Item(21,12).Value = "Segunda"
$ws.Cells.Item(21,13).Value = 10
$ws.Cells.Item(21,14).Value = 330000
$ws.Cells.Item(21,15).Value = 330000
$ws.Cells.Item(21,16).Value = 330000
$ws.Cells.Item(21,17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(21,18).Value = "Región Metropolitana"
$ws.Cells.Item(21,19).Value = 943
$ws.Cells.Item(21,20).Value = 350
